$d = $word.ActiveDocument

# 1) Replace the "I plan to use one of the following products..." sentence
#    with the new wording (formatting of the surrounding run is preserved
#    automatically because Find/Replace inherits the formatting of the
#    found text).
$old = "I plan to use one of the following products to separate saddle brackets and joist hangers from faces of piers and faces of beams and joists within 1 inch of the saddle brackets or joist hangers."
$new = "I plan to use a product like the following to separate saddle brackets from piers, to separate joist hangers from faces of beams and joists, and to protect beams and joists from the elements."
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null

# 2) Remove the five extra "here is another option" product hyperlinks
#    (Amazon POLYGUARD / Cylinnda / 3M tape / ELK / Stuck Tape) together with
#    the blank paragraphs that separated them, while keeping the first
#    (Lowes) link and the single trailing blank paragraph before the next
#    section ("Siding will extend..."). Locate them by their paragraph
#    text so the edit is resilient to any earlier shift in paragraph
#    indices.
$startPara = $null
$endPara = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -match "POLYGUARD-Poly-Wall-Decking-Adhering") {
        $startPara = $d.Paragraphs.Item($i - 1)
    }
    if ($t -match "Stuck-Tape-Multi-Purpose-Polyethylene-Protection") {
        $endPara = $p
    }
}
if ($startPara -ne $null -and $endPara -ne $null) {
    $rng = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $rng.Delete()
}

# 3) Remove the now-redundant blank paragraph that used to sit between the
#    "I plan to use..." sentence and the first (Lowes) product link.
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -match "^\s*$") {
        $next = $d.Paragraphs.Item($i + 1)
        $nt = $next.Range.Text
        if ($nt -match "lowes\.com/pd/APOC-Self-Bond") {
            $p.Range.Delete()
            break
        }
    }
}
